$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as Text so Excel does not
# reinterpret values like "1.00" or "66.865.12" as numbers/dates.
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "D9", "E9", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "E15", "D16", "E16", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "E24", "D25", "E25", "D26", "E26", "E27", "D28", "E28", "E29", "E30", "B31", "C31", "D31", "E31", "B32", "C32", "D32", "E32", "B33", "C33", "D33", "E33", "B34", "C34", "D34", "E34", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "B39", "C39", "D39", "E39", "B40", "C40", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "B45", "C45", "D45", "E45", "B46", "C46", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "E51")
foreach ($c in $cells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "66.865.12"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "3.515.04"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "584.04"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Value = "177.35"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").Value = "3.516.95"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D11").Value = "6.91"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("D12").Value = "0.422"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").Value = "4.120.91"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "30.62"
$ws.Range("E14").Value = "  -4.02%  "
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "66.886.90"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "3.502.62"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "6.11"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").Value = "14.06"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "380.94"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").Value = "7.86"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "0.536"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "71.58"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "9.91"
$ws.Range("E28").Value = "  -4.24%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "24.65"
$ws.Range("E31").Value = "  +4.90%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "6.00"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "2.02"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.37"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "7.19"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").Value = "1.56"
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("D38").Value = "158.46"
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "29.12"
$ws.Range("E39").Value = "  +11.39%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "0.889"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").Value = "1.80"
$ws.Range("E41").Value = "  -3.40%  "
$ws.Range("D42").Value = "2.64"
$ws.Range("E42").Value = "  -3.37%  "
$ws.Range("D43").Value = "6.60"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").Value = "4.53"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0707"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.718.29"
$ws.Range("E46").Value = "  -4.12%  "
$ws.Range("D47").Value = "25.64"
$ws.Range("E47").Value = "  -5.52%  "
$ws.Range("D48").Value = "40.50"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "0.0299"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "327.07"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("E51").Value = "  -1.92%  "
